$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for columns B, D, E, F, G, I, K across rows 2-25
# (updated line-loading results for the "380 kV" case)
$newValues = [ordered]@{
    2 = @{ "B" = 0.01842893888063202; "D" = 0.07997195081330233; "E" = 0.4080488715892159; "F" = 2.160681686102549; "G" = 0.002470119895837944; "I" = 0.8127539948200422; "K" = 1.474961274610394 }
    3 = @{ "B" = 0.0162820418710794; "D" = 0.07819488545445807; "E" = 0.3556212550602567; "F" = 2.063974498097792; "G" = 0.002476569284796371; "I" = 0.7977638552574149; "K" = 1.323869618261085 }
    4 = @{ "B" = 0.01496719093466936; "D" = 0.07713004808464063; "E" = 0.3235871875810261; "F" = 2.005791947699009; "G" = 0.002480726453144666; "I" = 0.7887752440649294; "K" = 1.232019525583439 }
    5 = @{ "B" = 0.01443239619419501; "D" = 0.07670255432707052; "E" = 0.3105689625498229; "F" = 1.982377654601379; "G" = 0.002482470334495334; "I" = 0.7851655365629995; "K" = 1.194815124894319 }
    6 = @{ "B" = 0.01434366054633784; "D" = 0.07663195297798353; "E" = 0.3084093646869661; "F" = 1.978507430268621; "G" = 0.002482762918867586; "I" = 0.7845693346549893; "K" = 1.188650786295 }
    7 = @{ "B" = 0.01495997415275951; "D" = 0.07712425692403002; "E" = 0.3234114783872144; "F" = 2.005474984412587; "G" = 0.002480749769916323; "I" = 0.7887263480192104; "K" = 1.231516870178268 }
    8 = @{ "B" = 0.01768809675392902; "D" = 0.07935365675033523; "E" = 0.3899371777388438; "F" = 2.127085632980851; "G" = 0.002472302839597137; "I" = 0.8075401187824909; "K" = 1.422669895991817 }
    9 = @{ "B" = 0.02305713625050743; "D" = 0.08394221359944254; "E" = 0.5218022583690072; "F" = 2.375303366918217; "G" = 0.00245729364122482; "I" = 0.8461872881795003; "K" = 1.805133970932388 }
    10 = @{ "B" = 0.02700391582787631; "D" = 0.08745719107415084; "E" = 0.6197821132993795; "F" = 2.563974407367397; "G" = 0.002447200928893857; "I" = 0.8757146855362521; "K" = 2.091232604743936 }
    11 = @{ "B" = 0.02879787928779365; "D" = 0.08909003739431398; "E" = 0.6646478438530323; "F" = 2.651257030991701; "G" = 0.00244280949881089; "I" = 0.889407722419719; "K" = 2.222598193319698 }
    12 = @{ "B" = 0.02947682593018186; "D" = 0.08971343087152661; "E" = 0.6816839433008823; "F" = 2.684524376585813; "G" = 0.002441175081994247; "I" = 0.8946315370181992; "K" = 2.272526176546762 }
    13 = @{ "B" = 0.02933062275976539; "D" = 0.0895789432177736; "E" = 0.6780127833470431; "F" = 2.677349984854345; "G" = 0.002441525817666091; "I" = 0.8935047633977291; "K" = 2.261765044607557 }
    14 = @{ "B" = 0.0288537454511868; "D" = 0.0891412215818832; "E" = 0.6660484597415177; "F" = 2.653989605310699; "G" = 0.002442674463975258; "I" = 0.8898367099889413; "K" = 2.226702093263896 }
    15 = @{ "B" = 0.02856158841015599; "D" = 0.08887377070608693; "E" = 0.6587261352453737; "F" = 2.639708919709051; "G" = 0.002443381751350948; "I" = 0.887594974062111; "K" = 2.205249041868854 }
    16 = @{ "B" = 0.02688663155195314; "D" = 0.08735117922358882; "E" = 0.6168563019884346; "F" = 2.558300079991682; "G" = 0.002447491921845158; "I" = 0.8748251519306933; "K" = 2.082672702911907 }
    17 = @{ "B" = 0.02585860069854817; "D" = 0.08642593380092478; "E" = 0.5912487656775625; "F" = 2.508735379000655; "G" = 0.002450064402774954; "I" = 0.8670588361397051; "K" = 2.007793332499546 }
    18 = @{ "B" = 0.02526717913245591; "D" = 0.08589693250874575; "E" = 0.5765474280154734; "F" = 2.480363677275591; "G" = 0.002451562844080452; "I" = 0.8626163483981344; "K" = 1.964838868454876 }
    19 = @{ "B" = 0.02506691792333271; "D" = 0.08571836109283026; "E" = 0.5715743967279963; "F" = 2.47078078584687; "G" = 0.002452073428644727; "I" = 0.8611163731923668; "K" = 1.950314606694349 }
    20 = @{ "B" = 0.02596805062689356; "D" = 0.08652409767562119; "E" = 0.5939718621996377; "F" = 2.513997444653256; "G" = 0.002449788611757633; "I" = 0.8678830313844799; "K" = 2.01575249911474 }
    21 = @{ "B" = 0.02899382785059146; "D" = 0.0892696518253473; "E" = 0.6695613784808643; "F" = 2.660845222765033; "G" = 0.002442336306481171; "I" = 0.8909130513293064; "K" = 2.236995912771704 }
    22 = @{ "B" = 0.03096900120673496; "D" = 0.09109367747649344; "E" = 0.7192360536668048; "F" = 2.758076622350558; "G" = 0.002437631946538008; "I" = 0.906189743929886; "K" = 2.382659755445559 }
    23 = @{ "B" = 0.02991508816180044; "D" = 0.09011738088133825; "E" = 0.6926974207504628; "F" = 2.706065266965197; "G" = 0.002440127617781839; "I" = 0.8980153273790421; "K" = 2.304815906631575 }
    24 = @{ "B" = 0.02591856952452076; "D" = 0.08647970866188359; "E" = 0.5927406858592406; "F" = 2.511618079914001; "G" = 0.00244991323593776; "I" = 0.8675103428021913; "K" = 2.012153865538892 }
    25 = @{ "B" = 0.02160360523298266; "D" = 0.08267649308824332; "E" = 0.4859539798657124; "F" = 2.307073544936117; "G" = 0.002461188930882586; "I" = 0.8355378675724481; "K" = 1.700805229751381 }
}

foreach ($row in $newValues.Keys) {
    $rowData = $newValues[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}

Write-Output "Updated $($newValues.Count) rows (380 kV case)"
